{"js": "// Map of old \"NNN\u00d7N=\" expressions to their new replacements, as described\n// by the diff. Each key is unique within the document, so a straightforward\n// search-and-replace (scoped to exact, case-sensitive matches) reproduces\n// the change precisely.\nconst replacements = [\n  [\"169\u00d76=\", \"950\u00d77=\"],\n  [\"228\u00d77=\", \"647\u00d73=\"],\n  [\"162\u00d72=\", \"862\u00d75=\"],\n  [\"264\u00d73=\", \"900\u00d76=\"],\n  [\"319\u00d77=\", \"417\u00d75=\"],\n  [\"688\u00d74=\", \"975\u00d74=\"],\n  [\"786\u00d79=\", \"469\u00d76=\"],\n  [\"794\u00d79=\", \"225\u00d78=\"],\n  [\"486\u00d72=\", \"889\u00d78=\"],\n  [\"875\u00d72=\", \"193\u00d76=\"],\n  [\"743\u00d72=\", \"407\u00d77=\"],\n  [\"214\u00d76=\", \"666\u00d75=\"],\n  [\"811\u00d78=\", \"458\u00d78=\"],\n  [\"824\u00d75=\", \"449\u00d77=\"],\n  [\"158\u00d77=\", \"641\u00d74=\"],\n  [\"164\u00d72=\", \"775\u00d78=\"],\n  [\"276\u00d78=\", \"603\u00d72=\"],\n  [\"779\u00d73=\", \"285\u00d79=\"],\n  [\"728\u00d75=\", \"177\u00d79=\"],\n  [\"358\u00d72=\", \"848\u00d72=\"],\n  [\"397\u00d75=\", \"850\u00d79=\"],\n  [\"162\u00d79=\", \"212\u00d72=\"],\n  [\"106\u00d77=\", \"212\u00d75=\"],\n  [\"367\u00d76=\", \"326\u00d79=\"],\n  [\"711\u00d72=\", \"441\u00d76=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each \"NNN\u00d7N=\" expression in the document with its updated value,\n# as described by the diff. Each \"before\" string is unique in the document,\n# so Find/Replace (scoped to the exact text, whole document) reproduces the\n# change precisely without touching anything else.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"169\u00d76=\", \"950\u00d77=\"),\n    @(\"228\u00d77=\", \"647\u00d73=\"),\n    @(\"162\u00d72=\", \"862\u00d75=\"),\n    @(\"264\u00d73=\", \"900\u00d76=\"),\n    @(\"319\u00d77=\", \"417\u00d75=\"),\n    @(\"688\u00d74=\", \"975\u00d74=\"),\n    @(\"786\u00d79=\", \"469\u00d76=\"),\n    @(\"794\u00d79=\", \"225\u00d78=\"),\n    @(\"486\u00d72=\", \"889\u00d78=\"),\n    @(\"875\u00d72=\", \"193\u00d76=\"),\n    @(\"743\u00d72=\", \"407\u00d77=\"),\n    @(\"214\u00d76=\", \"666\u00d75=\"),\n    @(\"811\u00d78=\", \"458\u00d78=\"),\n    @(\"824\u00d75=\", \"449\u00d77=\"),\n    @(\"158\u00d77=\", \"641\u00d74=\"),\n    @(\"164\u00d72=\", \"775\u00d78=\"),\n    @(\"276\u00d78=\", \"603\u00d72=\"),\n    @(\"779\u00d73=\", \"285\u00d79=\"),\n    @(\"728\u00d75=\", \"177\u00d79=\"),\n    @(\"358\u00d72=\", \"848\u00d72=\"),\n    @(\"397\u00d75=\", \"850\u00d79=\"),\n    @(\"162\u00d79=\", \"212\u00d72=\"),\n    @(\"106\u00d77=\", \"212\u00d75=\"),\n    @(\"367\u00d76=\", \"326\u00d79=\"),\n    @(\"711\u00d72=\", \"441\u00d76=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
